$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text edits (partial in-place text replacement) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "34"

$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "8/21/2023"
$c9.Characters(47, 9).Text = "8/27/2023"

# --- Data cell edits ---
$ws.Range("F14").Value = "0"
$ws.Range("L15").Value = -21.428571428571
$ws.Range("N15").Value = -71.794871794871
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 8
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 106
$ws.Range("J16").Value = 72
$ws.Range("K16").Value = 47.222222222222
$ws.Range("L16").Value = 92.727272727272
$ws.Range("M16").Value = -42.702702702702
$ws.Range("N16").Value = -83.489096573208
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -44.444444444444
$ws.Range("F17").Value = 38
$ws.Range("H17").Value = 22.58064516129
$ws.Range("I17").Value = 302
$ws.Range("J17").Value = 263
$ws.Range("K17").Value = 14.828897338403
$ws.Range("L17").Value = 66.850828729281
$ws.Range("M17").Value = 50.248756218905
$ws.Range("N17").Value = -41.586073500967
$ws.Range("C18").Value = "0"
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 25
$ws.Range("J18").Value = 59
$ws.Range("K18").Value = 55.93220338983
$ws.Range("L18").Value = 119.047619047619
$ws.Range("M18").Value = -51.322751322751
$ws.Range("N18").Value = -92.103004291845
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 14.285714285714
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 27.586206896551
$ws.Range("I19").Value = 271
$ws.Range("J19").Value = 231
$ws.Range("K19").Value = 17.316017316017
$ws.Range("L19").Value = 79.470198675496
$ws.Range("M19").Value = 1.119402985074
$ws.Range("N19").Value = -27.150537634408
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 86.666666666666
$ws.Range("I20").Value = 102
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = 45.714285714285
$ws.Range("L20").Value = 131.818181818182
$ws.Range("M20").Value = -23.880597014925
$ws.Range("N20").Value = -88.577827547592
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -26.923076923076
$ws.Range("F21").Value = 123
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = 23
$ws.Range("I21").Value = 894
$ws.Range("J21").Value = 706
$ws.Range("K21").Value = 26.628895184136
$ws.Range("L21").Value = 81.70731707317
$ws.Range("M21").Value = -10.956175298804
$ws.Range("N21").Value = -75.43956043956
$ws.Range("C23").Value = "0"
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 38
$ws.Range("K23").Value = 63.157894736842
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -9.677419354838
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = -10.169491525423
$ws.Range("I24").Value = 818
$ws.Range("J24").Value = 780
$ws.Range("K24").Value = 4.871794871794
$ws.Range("L24").Value = 57.307692307692
$ws.Range("M24").Value = -17.206477732793
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 8.333333333333
$ws.Range("F25").Value = 63
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = 43.181818181818
$ws.Range("I25").Value = 479
$ws.Range("J25").Value = 422
$ws.Range("K25").Value = 13.507109004739
$ws.Range("L25").Value = 38.840579710144
$ws.Range("M25").Value = -40.496894409937
$ws.Range("C26").Value = "0"
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -66.666666666666
$ws.Range("I26").Value = 17
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = -15
$ws.Range("L26").Value = -29.166666666666
$ws.Range("C27").Value = 7
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 600
$ws.Range("F27").Value = 12
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 63
$ws.Range("J27").Value = 48
$ws.Range("K27").Value = 31.25
$ws.Range("L27").Value = 103.225806451613
$ws.Range("F28").Value = 2
$ws.Range("L28").Value = -13.636363636363
$ws.Range("F29").Value = 2
$ws.Range("L29").Value = -15
$ws.Range("F30").Value = "0"
